$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.02871572971344
$ws.Range("B1").Value = 3.467015743255615
$ws.Range("C1").Value = 2.66535210609436
$ws.Range("D1").Value = 2.451266765594482
$ws.Range("E1").Value = 1.999853134155273
